# Applies the diff:
#  1. Moves the lone "_GoBack" bookmark from the end of the
#     "...координаты черепахи" paragraph into the middle of the
#     "Используется символы Unicode(50)" paragraph (right after the new
#     leading "!").
#  2. Rewrites that paragraph's runs from
#        "Используется символы " + "Unicode" + "(" + "50)"
#     to
#        "!" + <bookmark> + "Используется" + " символы " + "Unicode" + "(50)"
#     while preserving each run's original character formatting
#     (rFonts/sz/szCs/lang) exactly.
#
# Implementation notes (quirks of this COM-interop runtime worked around
# here):
#   * Range objects do not "live track" - after any edit that changes
#     text length, previously captured Start/End offsets of other Range
#     objects are stale, so ranges are re-derived (re-fetch the paragraph,
#     re-run Find) after every edit rather than reused.
#   * A Range's FormattedText getter is resolved lazily against the
#     *current* document at its original absolute offsets rather than
#     being a detached snapshot, so the clone-a-formatted-run trick
#     (capture FormattedText, then assign it elsewhere to create a new
#     run with identical rPr) only works if the paste happens immediately
#     after the capture, with no other edit in between. So both brand new
#     runs ("Unicode" and the run that becomes "(50)") are cloned first,
#     while the paragraph is still untouched, and every other edit
#     (retexting runs, moving the bookmark, deleting the old lead run)
#     happens afterwards.
#   * A paragraph Range's .End is *after* its trailing paragraph mark, so
#     "append inside this paragraph" inserts at `.End - 1`.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive current text.
$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("Используется символы") -and $t.Contains("Unicode")) {
        $paraIndex = $i
        break
    }
}

function Get-ParaRange {
    $p = $d.Paragraphs.Item($paraIndex)
    return $d.Range($p.Range.Start, $p.Range.End)
}

# --- Clone "Unicode" (rFonts/sz/szCs/lang=en-US) to a new run appended
#     at the end of the paragraph. Must paste immediately after capture.

$scope = Get-ParaRange
$rUnicode = $d.Range($scope.Start, $scope.End)
$rUnicode.Find.Execute("Unicode") | Out-Null
$ftUnicode = $rUnicode.FormattedText
$scope = Get-ParaRange
$insertAt = $scope.End - 1
$d.Range($insertAt, $insertAt).FormattedText = $ftUnicode

# --- Clone "50)" (rFonts/sz/szCs, no lang) to a new run appended right
#     after that. Will be retexted to "(50)" later. Must paste
#     immediately after capture.

$scope = Get-ParaRange
$rFifty = $d.Range($scope.Start, $scope.End)
$rFifty.Find.Execute("50)") | Out-Null
$ftPlain = $rFifty.FormattedText
$scope = Get-ParaRange
$insertAt = $scope.End - 1
$d.Range($insertAt, $insertAt).FormattedText = $ftPlain

# --- Now do the in-place edits on the original four runs. ---

# R2 "Unicode" -> "!"  (first occurrence of "Unicode" left in the
# paragraph - the clone appended above sits after the still-unedited
# "(" / "50)" runs, so it is unaffected.)
$scope = Get-ParaRange
$rUnicode2 = $d.Range($scope.Start, $scope.End)
$rUnicode2.Find.Execute("Unicode") | Out-Null
$rUnicode2.Text = "!"

# Bookmark "_GoBack" goes right after the new "!" run. Bookmarks.Add moves
# an existing bookmark of the same name rather than duplicating it, so
# this both removes "_GoBack" from its old location (end of the
# "...координаты черепахи" paragraph) and places it here.
$scope = Get-ParaRange
$rBang = $d.Range($scope.Start, $scope.End)
$rBang.Find.Execute("!") | Out-Null
$bmPoint = $d.Range($rBang.End, $rBang.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# R3 "(" -> "Используется"
$scope = Get-ParaRange
$rOpenParen = $d.Range($scope.Start, $scope.End)
$rOpenParen.Find.Execute("(") | Out-Null
$rOpenParen.Text = "Используется"

# R4 "50)" -> " символы "
$scope = Get-ParaRange
$rFiftyOrig = $d.Range($scope.Start, $scope.End)
$rFiftyOrig.Find.Execute("50)") | Out-Null
$rFiftyOrig.Text = " символы "

# Cloned "50)" run (appended earlier, right after the cloned "Unicode")
# -> "(50)"
$scope = Get-ParaRange
$rNew50 = $d.Range($scope.Start, $scope.End)
$rNew50.Find.Execute("50)") | Out-Null
$rNew50.Text = "(50)"

# --- Remove the now-obsolete leading run "Используется символы " ---
$scope = Get-ParaRange
$rLead = $d.Range($scope.Start, $scope.End)
$rLead.Find.Execute("Используется символы ") | Out-Null
$rLead.Text = ""
